$wb = $excel.ActiveWorkbook

$wsWebsites = $wb.Worksheets.Item("Useful Websites")
$wsLearning = $wb.Worksheets.Item("Learning Path")

# "Useful Websites" sheet gets the existing URL string in A1
$wsWebsites.Range("A1").Value = "https://www.datacamp.com/community/tutorials/machine-learning-python"

# "Learning Path" sheet gets two new rows, replacing the URL that used to be there
$wsLearning.Range("A1").Value = "All Harward lectures CS109"
$wsLearning.Range("A2").Value = "For practical strated with Iris dataset"

# Update the active selection on the "Learning Path" sheet to C7
$wsLearning.Activate()
$wsLearning.Range("C7").Select()
